$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert 3 new rows before the current row 10 (new negative-scenario
#    test case TC_EC_0005 ends up at the old row 10 position, i.e. row 13)
# ---------------------------------------------------------------------
$ws.Rows("10:12").Insert()
$ws.Rows("10:12").RowHeight = 15

# ---------------------------------------------------------------------
# 2. Populate the newly inserted rows with the new negative test case
#    (TC_EC_0004 - getpet negative scenario)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "TC_EC_0004"

$ws.Range("B10").Value = "comment"
$ws.Range("C10").Value = "petid"
$ws.Range("D10").Value = "status code"
$ws.Range("E10").Value = "message"
$ws.Range("F10").Value = "type"
$ws.Range("G10").Value = "code"
$ws.Range("B1").Copy()
$ws.Range("B10:G10").PasteSpecial(-4122)

$ws.Range("B11").Value = "getpet"

$ws.Range("C11").Value = "'48956"
$ws.Range("C8").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("D11").Value = "'404"
$ws.Range("D2").Copy()
$ws.Range("D11").PasteSpecial(-4122)

$ws.Range("E11").Value = "Pet not found"
$ws.Range("F11").Value = "error"

$ws.Range("G11").Value = "'1"
$ws.Range("D2").Copy()
$ws.Range("G11").PasteSpecial(-4122)

$ws.Range("A12").Value = "End"
$ws.Range("C12").Value = $null
$ws.Range("C9").Copy()
$ws.Range("C12").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. The old row 10 block (now shifted to row 13) becomes the new
#    negative-scenario test case id TC_EC_0005 instead of TC_EC_0004
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "TC_EC_0005"

# ---------------------------------------------------------------------
# 4. Header-style cleanup: the header rows (row 1 and row 3, and their
#    duplicated block at rows 13/15) now render all header cells with
#    the same font/alignment as column B (previously only column B used
#    it while C/D used a plain style).
# ---------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$ws.Range("B3:D3").PasteSpecial(-4122)
$ws.Range("C13:D13").PasteSpecial(-4122)
$ws.Range("B15:D15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. Restore the selection the author left active (B13:D13) and clear
#    the clipboard marching ants.
# ---------------------------------------------------------------------
$ws.Range("B13:D13").Select()
